# The commit swaps the deck's theme palette from the "Integral" design to
# the stock "Office Theme" palette (clrScheme dk2/lt2/accent1-6/hlink/
# folHlink). dk1/lt1 (black/white) and the font scheme (Arial everywhere)
# are identical between the two themes, so only the ten distinguishing
# colors need to change.
#
# This host has no filesystem-backed "Apply a theme file" operation
# (Theme.Save/ApplyTheme are sandboxed out) — theme edits go through the
# live ThemeColorScheme.Colors(i).RGB / ThemeFontScheme.MajorFont/MinorFont
# members instead, exactly like the in-app "Colors" / "Fonts" galleries.

function ToRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

# MsoThemeColorSchemeIndex order: 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5-10 accent1-6, 11 hlink, 12 folHlink.
$officeTheme = @{
    3  = "44546A"  # dk2
    4  = "E7E6E6"  # lt2
    5  = "5B9BD5"  # accent1
    6  = "ED7D31"  # accent2
    7  = "A5A5A5"  # accent3
    8  = "FFC000"  # accent4
    9  = "4472C4"  # accent5
    10 = "70AD47"  # accent6
    11 = "0563C1"  # hlink
    12 = "954F72"  # folHlink
}

foreach ($idx in $officeTheme.Keys) {
    $colors.Colors($idx).RGB = ToRGB $officeTheme[$idx]
}
